$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.82499266666667
$ws.Range("H2").Value = 56.474978
$ws.Range("I2").Value = 0.06886869772378311
$ws.Range("J2").Value = 0.0688686977237831
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 519.8509316666667
$ws.Range("N2").Value = 1559.552795
$ws.Range("O2").Value = 0.7998736652701117
$ws.Range("P2").Value = 0.7998736652701117
$ws.Range("Q2").Value = 9786.189976384836
$ws.Range("R2").Value = 88075.70978746351
$ws.Range("S2").Value = 0.0550862576707018
$ws.Range("T2").Value = 0.05508625767070179

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.82499266666667
$ws.Range("H3").Value = 56.474978
$ws.Range("I3").Value = 0.06886869772378311
$ws.Range("J3").Value = 0.0688686977237831
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.259691
$ws.Range("N3").Value = 3.779073
$ws.Range("O3").Value = 0.001938235743941786
$ws.Range("P3").Value = 0.001938235743941786
$ws.Range("Q3").Value = 23.713673837266
$ws.Range("R3").Value = 213.423064535394
$ws.Range("S3").Value = 0.0001334837715669587
$ws.Range("T3").Value = 0.0001334837715669587

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.82499266666667
$ws.Range("H4").Value = 56.474978
$ws.Range("I4").Value = 0.06886869772378311
$ws.Range("J4").Value = 0.0688686977237831
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 75.63123766666666
$ws.Range("N4").Value = 226.893713
$ws.Range("O4").Value = 0.1163707355248944
$ws.Range("P4").Value = 0.1163707355248944
$ws.Range("Q4").Value = 1423.757494445924
$ws.Range("R4").Value = 12813.81745001331
$ws.Range("S4").Value = 0.008014301008758258
$ws.Range("T4").Value = 0.008014301008758256

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.82499266666667
$ws.Range("H5").Value = 56.474978
$ws.Range("I5").Value = 0.06886869772378311
$ws.Range("J5").Value = 0.0688686977237831
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 53.174438
$ws.Range("N5").Value = 159.523314
$ws.Range("O5").Value = 0.08181736346105226
$ws.Range("P5").Value = 0.08181736346105226
$ws.Range("Q5").Value = 1001.008405404121
$ws.Range("R5").Value = 9009.075648637092
$ws.Range("S5").Value = 0.005634655272756105
$ws.Range("T5").Value = 0.005634655272756104

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 121.8208923333333
$ws.Range("H6").Value = 365.462677
$ws.Range("I6").Value = 0.4456653109566078
$ws.Range("J6").Value = 0.4456653109566078
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 519.8509316666667
$ws.Range("N6").Value = 1559.552795
$ws.Range("O6").Value = 0.7998736652701117
$ws.Range("P6").Value = 0.7998736652701117
$ws.Range("Q6").Value = 63328.70437594804
$ws.Range("R6").Value = 569958.3393835323
$ws.Range("S6").Value = 0.356475945758606
$ws.Range("T6").Value = 0.3564759457586059

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 121.8208923333333
$ws.Range("H7").Value = 365.462677
$ws.Range("I7").Value = 0.4456653109566078
$ws.Range("J7").Value = 0.4456653109566078
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.259691
$ws.Range("N7").Value = 3.779073
$ws.Range("O7").Value = 0.001938235743941786
$ws.Range("P7").Value = 0.001938235743941786
$ws.Range("Q7").Value = 153.456681684269
$ws.Range("R7").Value = 1381.110135158421
$ws.Range("S7").Value = 0.000863804435531028
$ws.Range("T7").Value = 0.000863804435531028

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 121.8208923333333
$ws.Range("H8").Value = 365.462677
$ws.Range("I8").Value = 0.4456653109566078
$ws.Range("J8").Value = 0.4456653109566078
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 75.63123766666666
$ws.Range("N8").Value = 226.893713
$ws.Range("O8").Value = 0.1163707355248944
$ws.Range("P8").Value = 0.1163707355248944
$ws.Range("Q8").Value = 9213.464860827746
$ws.Range("R8").Value = 82921.1837474497
$ws.Range("S8").Value = 0.05186240003395121
$ws.Range("T8").Value = 0.05186240003395121

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 121.8208923333333
$ws.Range("H9").Value = 365.462677
$ws.Range("I9").Value = 0.4456653109566078
$ws.Range("J9").Value = 0.4456653109566078
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 53.174438
$ws.Range("N9").Value = 159.523314
$ws.Range("O9").Value = 0.08181736346105226
$ws.Range("P9").Value = 0.08181736346105226
$ws.Range("Q9").Value = 6477.75748648351
$ws.Range("R9").Value = 58299.81737835158
$ws.Range("S9").Value = 0.03646316072851966
$ws.Range("T9").Value = 0.03646316072851966

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 87.673585
$ws.Range("H10").Value = 263.020755
$ws.Range("I10").Value = 0.3207419907481189
$ws.Range("J10").Value = 0.3207419907481188
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 519.8509316666667
$ws.Range("N10").Value = 1559.552795
$ws.Range("O10").Value = 0.7998736652701117
$ws.Range("P10").Value = 0.7998736652701117
$ws.Range("Q10").Value = 45577.1948448067
$ws.Range("R10").Value = 410194.7536032603
$ws.Range("S10").Value = 0.2565530717457301
$ws.Range("T10").Value = 0.2565530717457301

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 87.673585
$ws.Range("H11").Value = 263.020755
$ws.Range("I11").Value = 0.3207419907481189
$ws.Range("J11").Value = 0.3207419907481188
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.259691
$ws.Range("N11").Value = 3.779073
$ws.Range("O11").Value = 0.001938235743941786
$ws.Range("P11").Value = 0.001938235743941786
$ws.Range("Q11").Value = 110.441625962235
$ws.Range("R11").Value = 993.974633660115
$ws.Range("S11").Value = 0.0006216735910510496
$ws.Range("T11").Value = 0.0006216735910510495

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 87.673585
$ws.Range("H12").Value = 263.020755
$ws.Range("I12").Value = 0.3207419907481189
$ws.Range("J12").Value = 0.3207419907481188
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 75.63123766666666
$ws.Range("N12").Value = 226.893713
$ws.Range("O12").Value = 0.1163707355248944
$ws.Range("P12").Value = 0.1163707355248944
$ws.Range("Q12").Value = 6630.861744223702
$ws.Range("R12").Value = 59677.75569801332
$ws.Range("S12").Value = 0.03732498137707745
$ws.Range("T12").Value = 0.03732498137707745

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 87.673585
$ws.Range("H13").Value = 263.020755
$ws.Range("I13").Value = 0.3207419907481189
$ws.Range("J13").Value = 0.3207419907481188
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 53.174438
$ws.Range("N13").Value = 159.523314
$ws.Range("O13").Value = 0.08181736346105226
$ws.Range("P13").Value = 0.08181736346105226
$ws.Range("Q13").Value = 4661.993609820231
$ws.Range("R13").Value = 41957.94248838207
$ws.Range("S13").Value = 0.0262422640342603
$ws.Range("T13").Value = 0.0262422640342603

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 45.02666966666666
$ws.Range("H14").Value = 135.080009
$ws.Range("I14").Value = 0.1647240005714903
$ws.Range("J14").Value = 0.1647240005714903
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 519.8509316666667
$ws.Range("N14").Value = 1559.552795
$ws.Range("O14").Value = 0.7998736652701117
$ws.Range("P14").Value = 0.7998736652701117
$ws.Range("Q14").Value = 23407.15617606391
$ws.Range("R14").Value = 210664.4055845751
$ws.Range("S14").Value = 0.1317583900950739
$ws.Range("T14").Value = 0.1317583900950739

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 45.02666966666666
$ws.Range("H15").Value = 135.080009
$ws.Range("I15").Value = 0.1647240005714903
$ws.Range("J15").Value = 0.1647240005714903
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.259691
$ws.Range("N15").Value = 3.779073
$ws.Range("O15").Value = 0.001938235743941786
$ws.Range("P15").Value = 0.001938235743941786
$ws.Range("Q15").Value = 56.71969053907299
$ws.Range("R15").Value = 510.4772148516569
$ws.Range("S15").Value = 0.0003192739457927496
$ws.Range("T15").Value = 0.0003192739457927496

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 45.02666966666666
$ws.Range("H16").Value = 135.080009
$ws.Range("I16").Value = 0.1647240005714903
$ws.Range("J16").Value = 0.1647240005714903
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 75.63123766666666
$ws.Range("N16").Value = 226.893713
$ws.Range("O16").Value = 0.1163707355248944
$ws.Range("P16").Value = 0.1163707355248944
$ws.Range("Q16").Value = 3405.422754898157
$ws.Range("R16").Value = 30648.80479408341
$ws.Range("S16").Value = 0.01916905310510744
$ws.Range("T16").Value = 0.01916905310510744

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 45.02666966666666
$ws.Range("H17").Value = 135.080009
$ws.Range("I17").Value = 0.1647240005714903
$ws.Range("J17").Value = 0.1647240005714903
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 53.174438
$ws.Range("N17").Value = 159.523314
$ws.Range("O17").Value = 0.08181736346105226
$ws.Range("P17").Value = 0.08181736346105226
$ws.Range("Q17").Value = 2394.267854536647
$ws.Range("R17").Value = 21548.41069082982
$ws.Range("S17").Value = 0.0134772834255162
$ws.Range("T17").Value = 0.0134772834255162
